# Adding realtime and prediction function results to the observations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header (row 35), styled like the other "using loader vX" headers ---
$ws.Range("A35").Value = "The best parameters are {'C': 10.0, 'gamma': 1e-09} with a score of 0.98"
$ws.Range("A35").Style = "Calculation"

# --- New experiment rows (37-39) ---

# Row 37
$ws.Range("A37").Value = 13
$ws.Range("C37").Value = 10
$ws.Range("E37").Value = 0.000000001
$ws.Range("E37").NumberFormat = "0.00E+00"
$ws.Range("G37").Value = "linear"
$ws.Range("I37").Value = 0.90419161676646698
$ws.Range("K37").Value = "TTP+TNA+angels between fingers"
$ws.Range("L37").Value = 28

# Row 38
$ws.Range("A38").Value = 14
$ws.Range("C38").Value = 1000
$ws.Range("E38").Value = 0.01
$ws.Range("G38").Value = "linear"
$ws.Range("I38").Value = 0.90269461077844304
$ws.Range("K38").Value = "TTP+TNA+angels between fingers"
$ws.Range("L38").Value = 28

# Row 39
$ws.Range("A39").Value = 15
$ws.Range("C39").Value = 1000000
$ws.Range("E39").Value = 0.01
$ws.Range("G39").Value = "linear"
$ws.Range("I39").Value = 0.90269461077844304
$ws.Range("K39").Value = "TTP+TNA+angels between fingers"
$ws.Range("L39").Value = 28

# --- Update the view: scroll position and active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M43").Select()
